$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 366 (this pushes the existing rows 366-371
# down to become rows 369-374, carrying their data/styles with them).
$ws.Range("A366:A368").EntireRow.Insert()

# Fill the 3 newly-inserted rows (366-368) with the new weekly data.
$rowsData = @(
    @{ Row = 366; L = "1a amarillo"; M = 240; N = 3800; O = 4000; P = 3900; R = "Región de O'Higgins"; S = 244 },
    @{ Row = 367; L = "2a amarillo"; M = 240; N = 3000; O = 3500; P = 3250; R = "Región de O'Higgins"; S = 203 },
    @{ Row = 368; L = "3a amarillo"; M = 200; N = 2600; O = 2800; P = 2700; R = "Región de O'Higgins"; S = 169 }
)

foreach ($rd in $rowsData) {
    $r = $rd.Row
    $ws.Cells.Item($r, 1).Value = 7
    $ws.Cells.Item($r, 2).Value = "Terminal Hortofrutícola Agro Chillán"
    $ws.Cells.Item($r, 3).Value = "Ñuble"
    $ws.Cells.Item($r, 4).Value = 44448
    $ws.Cells.Item($r, 5).Value = 16
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100102
    $ws.Cells.Item($r, 8).Value = "Cítricos"
    $ws.Cells.Item($r, 9).Value = 100102003
    $ws.Cells.Item($r, 10).Value = "Limón"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = $rd.L
    $ws.Cells.Item($r, 13).Value = $rd.M
    $ws.Cells.Item($r, 14).Value = $rd.N
    $ws.Cells.Item($r, 15).Value = $rd.O
    $ws.Cells.Item($r, 16).Value = $rd.P
    $ws.Cells.Item($r, 17).Value = "`$/malla 16 kilos"
    $ws.Cells.Item($r, 18).Value = $rd.R
    $ws.Cells.Item($r, 19).Value = $rd.S
    $ws.Cells.Item($r, 20).Value = 16
}
